$wb = $excel.ActiveWorkbook
$wsView = $wb.Worksheets.Item("view")
$wsCtrl = $wb.Worksheets.Item("controleur")

# --- "view" sheet: mark D3 / D4 / D6 as done ("V") -----------------------
# D3 and D4 already use the "done" style (s=18); just stamp the value.
$wsView.Range("D3").Value = "V"
$wsView.Range("D4").Value = "V"

# D6 currently uses a different style (s=17); copy the formatting that D5
# already has (s=18, the "done" look) onto D6, then stamp its value too.
$wsView.Range("D5").Copy()
$wsView.Range("D6").PasteSpecial(-4122)   # xlPasteFormats
$wsView.Range("D6").Value = "V"

# Applying the bigger "done" font to D6 grows row 6's height, same as the
# other "done" rows (2, 3, 4, 5) which all sit at 29pt.
$wsView.Rows.Item(6).RowHeight = 29

# --- sharedStrings text tweak on "controleur" -----------------------------
$wsCtrl.Range("C16").Value = "`$listeUsers(objet utilisateur)"

# --- selections / active sheet -------------------------------------------
# Move the cursor on "view" without leaving it the active tab.
$wsView.Activate()
$wsView.Range("C7").Select()

# Re-activate "controleur" (the tab that should stay selected) and move its
# cursor too.
$wsCtrl.Activate()
$wsCtrl.Range("A2").Select()
